# Auto-generated Excel COM-interop script
# Applies numeric value corrections across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# per the scheduled-runner recomputation described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1153.921
$ws.Range("I28").Value = 1153.0625
$ws.Range("J28").Value = 1158.5
$ws.Range("K28").Value = 1153.0625
$ws.Range("L28").Value = 1158.5
$ws.Range("M28").Value = -668.0625
$ws.Range("N28").Value = -2128.5
$ws.Range("H41").Value = 1079.069
$ws.Range("I41").Value = 940.43475
$ws.Range("J41").Value = 1610.5
$ws.Range("K41").Value = 940.43475
$ws.Range("L41").Value = 1610.5
$ws.Range("M41").Value = -500.43475
$ws.Range("N41").Value = -2490.5
$ws.Range("H94").Value = 1792.1818
$ws.Range("I94").Value = 1792.1818
$ws.Range("K94").Value = 1792.1818
$ws.Range("M94").Value = -1341.1818
$ws.Range("H125").Value = 1910.3
$ws.Range("I125").Value = 769.0833
$ws.Range("K125").Value = 6921.7497
$ws.Range("M125").Value = -4461.7497
$ws.Range("H132").Value = 783.678
$ws.Range("I132").Value = 813.22
$ws.Range("K132").Value = 2439.66
$ws.Range("M132").Value = 90.34000000000015
$ws.Range("H138").Value = 2516.0505
$ws.Range("J138").Value = 2875.111
$ws.Range("L138").Value = 8625.332999999999
$ws.Range("N138").Value = -18905.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13543.871
$ws.Range("I32").Value = 8538.937
$ws.Range("K32").Value = 8538.937
$ws.Range("M32").Value = -8251.937
$ws.Range("H61").Value = 4447.84
$ws.Range("I61").Value = 2667
$ws.Range("K61").Value = 2667
$ws.Range("M61").Value = -2455
$ws.Range("H63").Value = 7966.3335
$ws.Range("J63").Value = 5699.5
$ws.Range("L63").Value = 5699.5
$ws.Range("N63").Value = -7071.5
$ws.Range("H66").Value = 7966.3335
$ws.Range("J66").Value = 5699.5
$ws.Range("L66").Value = 28497.5
$ws.Range("N66").Value = -35361.5
$ws.Range("H132").Value = 5378.1284
$ws.Range("I132").Value = 5694.727
$ws.Range("K132").Value = 17084.181
$ws.Range("M132").Value = -14554.181
$ws.Range("H133").Value = 72583.336
$ws.Range("J133").Value = 81571.42999999999
$ws.Range("L133").Value = 81571.42999999999
$ws.Range("N133").Value = -86631.42999999999
$ws.Range("H136").Value = 4447.84
$ws.Range("I136").Value = 2667
$ws.Range("K136").Value = 8001
$ws.Range("M136").Value = -5451

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 166666.67
$ws.Range("I35").Value = 100000
$ws.Range("K35").Value = 100000
$ws.Range("M35").Value = -99690
$ws.Range("H134").Value = 4446.84
$ws.Range("I134").Value = 3383.5881
$ws.Range("J134").Value = 6706.25
$ws.Range("K134").Value = 10150.7643
$ws.Range("L134").Value = 20118.75
$ws.Range("M134").Value = -7615.764299999999
$ws.Range("N134").Value = -25188.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1234.2858
$ws.Range("I22").Value = 1035
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1035
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -685
$ws.Range("N22").Value = -2200
$ws.Range("H105").Value = 706.2857
$ws.Range("I105").Value = 672.25
$ws.Range("J105").Value = 751.6667
$ws.Range("K105").Value = 672.25
$ws.Range("L105").Value = 751.6667
$ws.Range("M105").Value = 1074.75
$ws.Range("N105").Value = -4245.6667
$ws.Range("H107").Value = 4979.8335
$ws.Range("I107").Value = 833.5
$ws.Range("J107").Value = 17418.834
$ws.Range("K107").Value = 833.5
$ws.Range("L107").Value = 17418.834
$ws.Range("M107").Value = 1086.5
$ws.Range("N107").Value = -21258.834
$ws.Range("H125").Value = 65997.5
$ws.Range("J125").Value = 65997.5
$ws.Range("L125").Value = 65997.5
$ws.Range("N125").Value = -70917.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1456.7142
$ws.Range("I113").Value = 1567.1428
$ws.Range("J113").Value = 1401.5
$ws.Range("K113").Value = 4701.428400000001
$ws.Range("L113").Value = 4204.5
$ws.Range("M113").Value = -2531.428400000001
$ws.Range("N113").Value = -8544.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5026.1875
$ws.Range("I132").Value = 4535
$ws.Range("J132").Value = 6499.75
$ws.Range("K132").Value = 13605
$ws.Range("L132").Value = 19499.25
$ws.Range("M132").Value = -11075
$ws.Range("N132").Value = -24559.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1645.9333
$ws.Range("I82").Value = 911.7143
$ws.Range("J82").Value = 2288.375
$ws.Range("K82").Value = 911.7143
$ws.Range("L82").Value = 2288.375
$ws.Range("M82").Value = -550.7143
$ws.Range("N82").Value = -3010.375
$ws.Range("H85").Value = 1645.9333
$ws.Range("I85").Value = 911.7143
$ws.Range("J85").Value = 2288.375
$ws.Range("K85").Value = 911.7143
$ws.Range("L85").Value = 2288.375
$ws.Range("M85").Value = 336.2857
$ws.Range("N85").Value = -4784.375
$ws.Range("H93").Value = 10402272
$ws.Range("I93").Value = 2752.4211
$ws.Range("J93").Value = 43334084
$ws.Range("K93").Value = 2752.4211
$ws.Range("L93").Value = 43334084
$ws.Range("M93").Value = -1504.4211
$ws.Range("N93").Value = -43336580
$ws.Range("H97").Value = 18499.334
$ws.Range("J97").Value = 18499.334
$ws.Range("L97").Value = 18499.334
$ws.Range("N97").Value = -20481.334
$ws.Range("H98").Value = 31568
$ws.Range("J98").Value = 31568
$ws.Range("L98").Value = 31568
$ws.Range("N98").Value = -37558
$ws.Range("H99").Value = 15000
$ws.Range("I99").Value = 15000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 15000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -12005
$ws.Range("N99").ClearContents()
$ws.Range("H132").Value = 4509.943
$ws.Range("I132").Value = 3697.5386
$ws.Range("K132").Value = 11092.6158
$ws.Range("M132").Value = -8562.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18521556
$ws.Range("J81").Value = 41670628
$ws.Range("L81").Value = 83341256
$ws.Range("N81").Value = -83343378
$ws.Range("H84").Value = 18521556
$ws.Range("J84").Value = 41670628
$ws.Range("L84").Value = 416706280
$ws.Range("N84").Value = -416716888
$ws.Range("H122").Value = 3972.8462
$ws.Range("I122").Value = 3963.125
$ws.Range("K122").Value = 11889.375
$ws.Range("M122").Value = -9439.375
$ws.Range("H136").Value = 9618504
$ws.Range("I136").Value = 15627527
$ws.Range("K136").Value = 46882581
$ws.Range("M136").Value = -46880031

Write-Output "Applied all cell updates."